$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: English passage
$ws.Range("B2").Value = 'English'
$ws.Range("A2").Value = 'English is a West Germanic language first spoken in early medieval England which eventually became the leading language of international discourse in today''s world.[4][5][6] It is named after the Angles, one of the ancient Germanic peoples that migrated to the area of Great Britain that later took their name, England. Both names derive from Anglia, a peninsula on the Baltic Sea. English is most closely related to Frisian and Low Saxon, while its vocabulary has been significantly influenced by other Germanic languages, particularly Old Norse (a North Germanic language), as well as Latin and French.[7][8][9]. 
A chemical bond is a lasting attraction between atoms, ions or molecules that enables the formation of chemical compounds. The bond may result from the electrostatic force of attraction between oppositely charged ions as in ionic bonds or through the sharing of electrons as in covalent bonds. The strength of chemical bonds varies considerably; there are "strong bonds" or "primary bonds" such as covalent, ionic and metallic bonds, and "weak bonds" or "secondary bonds" such as dipole–dipole interactions, the London dispersion force and hydrogen bonding.'
$ws.Range("C2").Value = "How many moons does the earth have?"

# Row 3: Anglo-Saxon England passage
$ws.Range("A3").Value = 'Anglo-Saxon England was early medieval England, existing from the 5th to the 11th centuries from the end of Roman Britain until the Norman conquest in 1066. It consisted of various Anglo-Saxon kingdoms until 927 when it was united as the Kingdom of England by King Æthelstan (r. 927–939). It became part of the short-lived North Sea Empire of Cnut the Great, a personal union between England, Denmark and Norway in the 11th century.
The Anglo-Saxons were the members of Germanic-speaking groups who migrated to the southern half of the island of Great Britain from nearby northwestern Europe. Anglo-Saxon history thus begins during the period of sub-Roman Britain following the end of Roman control, and traces the establishment of Anglo-Saxon kingdoms in the 5th and 6th centuries (conventionally identified as seven main kingdoms: Northumbria, Mercia, East Anglia, Essex, Kent, Sussex, and Wessex), their Christianisation during the 7th century, the threat of Viking invasions and Danish settlers, the gradual unification of England under the Wessex hegemony during the 9th and 10th centuries, and ending with the Norman conquest of England by William the Conqueror in 1066. '
$ws.Range("B3").Value = "BIology,"
$ws.Range("C3").Value = 'When did the United States of America become a country?'
$ws.Range("D3").Value = 2002
$ws.Range("F3").Value = 1938
$ws.Range("H3").Value = 1900

$ws.Rows.Item(2).RowHeight = 15

$ws.Range("C5").Select()
